# Updates profit/price figures on several "Leve" rows across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR crafting sheets (scheduled market-price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 516.2105
$ws.Range("I33").Value = 275.15384
$ws.Range("K33").Value = 275.15384
$ws.Range("M33").Value = -46.15384

# Row 69
$ws.Range("H69").Value = 9565.478999999999
$ws.Range("I69").Value = 4999.5
$ws.Range("J69").Value = 10000.333
$ws.Range("K69").Value = 14998.5
$ws.Range("L69").Value = 30000.999
$ws.Range("M69").Value = -14124.5
$ws.Range("N69").Value = -31748.999

# Row 72
$ws.Range("H72").Value = 9565.478999999999
$ws.Range("I72").Value = 4999.5
$ws.Range("J72").Value = 10000.333
$ws.Range("K72").Value = 44995.5
$ws.Range("L72").Value = 90002.997
$ws.Range("M72").Value = -40627.5
$ws.Range("N72").Value = -98738.997

# Row 74
$ws.Range("H74").Value = 7319.1304
$ws.Range("I74").Value = 3466.5
$ws.Range("J74").Value = 9373.866
$ws.Range("K74").Value = 3466.5
$ws.Range("L74").Value = 9373.866
$ws.Range("M74").Value = -2530.5
$ws.Range("N74").Value = -11245.866

# Row 77
$ws.Range("H77").Value = 7319.1304
$ws.Range("I77").Value = 3466.5
$ws.Range("J77").Value = 9373.866
$ws.Range("K77").Value = 17332.5
$ws.Range("L77").Value = 46869.33
$ws.Range("M77").Value = -12652.5
$ws.Range("N77").Value = -56229.33

# Row 95
$ws.Range("H95").Value = 80415.5
$ws.Range("J95").Value = 80415.5
$ws.Range("L95").Value = 80415.5
$ws.Range("N95").Value = -85907.5

# Row 106
$ws.Range("H106").Value = 2848.8572
$ws.Range("I106").Value = 2849.7
$ws.Range("K106").Value = 2849.7
$ws.Range("M106").Value = -2218.7

# Row 113
$ws.Range("H113").Value = 3529.1667
$ws.Range("I113").Value = 3102.25
$ws.Range("J113").Value = 3742.625
$ws.Range("K113").Value = 3102.25
$ws.Range("L113").Value = 3742.625
$ws.Range("M113").Value = 151.75
$ws.Range("N113").Value = -10250.625

# Row 135
$ws.Range("H135").Value = 2651.7334
$ws.Range("I135").Value = 974.1667
$ws.Range("J135").Value = 3770.111
$ws.Range("K135").Value = 8767.5003
$ws.Range("L135").Value = 33930.999
$ws.Range("M135").Value = -6232.5003
$ws.Range("N135").Value = -39000.999

# Row 138
$ws.Range("H138").Value = 1833.5588
$ws.Range("I138").Value = 1159.174
$ws.Range("K138").Value = 3477.522
$ws.Range("M138").Value = 1662.478

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1668096.5
$ws.Range("I61").Value = 2223402.8
$ws.Range("J61").Value = 2178
$ws.Range("K61").Value = 2223402.8
$ws.Range("L61").Value = 2178
$ws.Range("M61").Value = -2223190.8
$ws.Range("N61").Value = -2602

# Row 74
$ws.Range("H74").Value = 4985
$ws.Range("I74").Value = 1973.75
$ws.Range("J74").Value = 9000
$ws.Range("K74").Value = 1973.75
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -1099.75
$ws.Range("N74").Value = -10748

# Row 77
$ws.Range("H77").Value = 4985
$ws.Range("I77").Value = 1973.75
$ws.Range("J77").Value = 9000
$ws.Range("K77").Value = 9868.75
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -5500.75
$ws.Range("N77").Value = -53736

# Row 122
$ws.Range("H122").Value = 1146.1428
$ws.Range("I122").Value = 1146.1428
$ws.Range("K122").Value = 3438.4284
$ws.Range("M122").Value = -988.4284000000002

# Row 132
$ws.Range("H132").Value = 793057.25
$ws.Range("I132").Value = 880619.3
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 2641857.9
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -2639327.9
$ws.Range("N132").Value = -20057

# Row 136
$ws.Range("H136").Value = 1668096.5
$ws.Range("I136").Value = 2223402.8
$ws.Range("J136").Value = 2178
$ws.Range("K136").Value = 6670208.399999999
$ws.Range("L136").Value = 6534
$ws.Range("M136").Value = -6667658.399999999
$ws.Range("N136").Value = -11634

$ws = $wb.Worksheets.Item("BSM")
# Row 132
$ws.Range("H132").Value = 47500
$ws.Range("J132").Value = 47500
$ws.Range("L132").Value = 47500
$ws.Range("N132").Value = -57620

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 800
$ws.Range("I16").Value = 800
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -513
$ws.Range("N16").ClearContents()

# Row 62
$ws.Range("H62").Value = 3117.5
$ws.Range("I62").Value = 5
$ws.Range("J62").Value = 3740
$ws.Range("K62").Value = 5
$ws.Range("L62").Value = 3740
$ws.Range("M62").Value = 619
$ws.Range("N62").Value = -4988

# Row 65
$ws.Range("H65").Value = 3117.5
$ws.Range("I65").Value = 5
$ws.Range("J65").Value = 3740
$ws.Range("K65").Value = 25
$ws.Range("L65").Value = 18700
$ws.Range("M65").Value = 3095
$ws.Range("N65").Value = -24940

# Row 97
$ws.Range("H97").Value = 26666.334
$ws.Range("J97").Value = 26666.334
$ws.Range("L97").Value = 26666.334
$ws.Range("N97").Value = -28648.334

# Row 113
$ws.Range("H113").Value = 800
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 800
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1370
$ws.Range("N113").ClearContents()

# Row 134
$ws.Range("H134").Value = 5899.9707
$ws.Range("J134").Value = 3010.25
$ws.Range("L134").Value = 9030.75
$ws.Range("N134").Value = -14100.75

$ws = $wb.Worksheets.Item("CUL")
# Row 47
$ws.Range("H47").Value = 1540
$ws.Range("I47").Value = 1175
$ws.Range("K47").Value = 3525
$ws.Range("M47").Value = -3094

# Row 68
$ws.Range("H68").Value = 6269.2856
$ws.Range("J68").Value = 7394.7646
$ws.Range("L68").Value = 22184.2938
$ws.Range("N68").Value = -23806.2938

# Row 71
$ws.Range("H71").Value = 6269.2856
$ws.Range("J71").Value = 7394.7646
$ws.Range("L71").Value = 66552.8814
$ws.Range("N71").Value = -74664.8814

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 323084.06
$ws.Range("I80").Value = 571425.9
$ws.Range("J80").Value = 3787.4285
$ws.Range("K80").Value = 571425.9
$ws.Range("L80").Value = 3787.4285
$ws.Range("M80").Value = -570427.9
$ws.Range("N80").Value = -5783.4285

# Row 83
$ws.Range("H83").Value = 323084.06
$ws.Range("I83").Value = 571425.9
$ws.Range("J83").Value = 3787.4285
$ws.Range("K83").Value = 2857129.5
$ws.Range("L83").Value = 18937.1425
$ws.Range("M83").Value = -2852137.5
$ws.Range("N83").Value = -28921.1425

# Row 97
$ws.Range("H97").Value = 1883.7273
$ws.Range("I97").Value = 1564.5264
$ws.Range("K97").Value = 1564.5264
$ws.Range("M97").Value = -1068.5264

# Row 122
$ws.Range("H122").Value = 65867.94
$ws.Range("I122").Value = 86379.586
$ws.Range("J122").Value = 16640
$ws.Range("K122").Value = 259138.758
$ws.Range("L122").Value = 49920
$ws.Range("M122").Value = -256688.758
$ws.Range("N122").Value = -54820

# Row 126
$ws.Range("H126").Value = 1391553.4
$ws.Range("I126").Value = 1854048.9
$ws.Range("J126").Value = 4066.6667
$ws.Range("K126").Value = 5562146.699999999
$ws.Range("L126").Value = 12200.0001
$ws.Range("M126").Value = -5559676.699999999
$ws.Range("N126").Value = -17140.0001

# Row 132
$ws.Range("H132").Value = 33740344
$ws.Range("I132").Value = 37488210
$ws.Range("J132").Value = 9564.333000000001
$ws.Range("K132").Value = 112464630
$ws.Range("L132").Value = 28692.999
$ws.Range("M132").Value = -112462100
$ws.Range("N132").Value = -33752.999

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3998
$ws.Range("I7").Value = 3730.8667
$ws.Range("J7").Value = 4999.75
$ws.Range("K7").Value = 3730.8667
$ws.Range("L7").Value = 4999.75
$ws.Range("M7").Value = -3618.8667
$ws.Range("N7").Value = -5223.75

# Row 61
$ws.Range("H61").Value = 2581.111
$ws.Range("I61").Value = 1941
$ws.Range("J61").Value = 4821.5
$ws.Range("K61").Value = 1941
$ws.Range("L61").Value = 4821.5
$ws.Range("M61").Value = -1739
$ws.Range("N61").Value = -5225.5

# Row 68
$ws.Range("H68").Value = 3640.6
$ws.Range("I68").Value = 4350.5
$ws.Range("J68").Value = 3167.3333
$ws.Range("K68").Value = 4350.5
$ws.Range("L68").Value = 3167.3333
$ws.Range("M68").Value = -3601.5
$ws.Range("N68").Value = -4665.3333

# Row 71
$ws.Range("H71").Value = 3640.6
$ws.Range("I71").Value = 4350.5
$ws.Range("J71").Value = 3167.3333
$ws.Range("K71").Value = 21752.5
$ws.Range("L71").Value = 15836.6665
$ws.Range("M71").Value = -18008.5
$ws.Range("N71").Value = -23324.6665

# Row 113
$ws.Range("H113").Value = 2581.111
$ws.Range("I113").Value = 1941
$ws.Range("J113").Value = 4821.5
$ws.Range("K113").Value = 1941
$ws.Range("L113").Value = 4821.5
$ws.Range("M113").Value = 229
$ws.Range("N113").Value = -9161.5

# Row 126
$ws.Range("H126").Value = 3998
$ws.Range("I126").Value = 3730.8667
$ws.Range("J126").Value = 4999.75
$ws.Range("K126").Value = 11192.6001
$ws.Range("L126").Value = 14999.25
$ws.Range("M126").Value = -8722.6001
$ws.Range("N126").Value = -19939.25

# Row 136
$ws.Range("H136").Value = 64413.7
$ws.Range("I136").Value = 2188.077
$ws.Range("J136").Value = 179975.58
$ws.Range("K136").Value = 6564.231000000001
$ws.Range("L136").Value = 539926.74
$ws.Range("M136").Value = -4014.231000000001
$ws.Range("N136").Value = -545026.74

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2545.7097
$ws.Range("I122").Value = 2073.08
$ws.Range("K122").Value = 6219.24
$ws.Range("M122").Value = -3769.24

# Row 126
$ws.Range("H126").Value = 6623.8
$ws.Range("I126").Value = 6487.091
$ws.Range("K126").Value = 19461.273
$ws.Range("M126").Value = -16991.273

# Row 132
$ws.Range("H132").Value = 4027827.5
$ws.Range("I132").Value = 4682474
$ws.Range("K132").Value = 14047422
$ws.Range("M132").Value = -14044892

# Row 141
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360
